# Generate Report for Handoff
# Updates the localization-status workbook: the "9c314e72-...md" file (row 3
# in every sheet) moves from "Handed back: in sync with en-US" to
# "Ready for handoff", and its per-language "Latest Handoff Datetime" is
# stamped with the new handoff timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-10 12:51:37"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-10 12:51:43"
